$d = $word.ActiveDocument
$apos = [char]0x2019

# ------------------------------------------------------------------
# 1) "First, I want to briefly go over ..." -> drop "briefly " and
#    leave Word's last-edit-position ("_GoBack") bookmark sitting
#    right before "go over", splitting that run in two.
# ------------------------------------------------------------------
$full = "First, I want to briefly go over our understanding of the assignment tasks."

$rng = $d.Content
$found = $rng.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find the intro sentence"
}
$sentenceStart = $rng.Start

$prefix = "First, I want to briefly "
$goOverOffset = $prefix.Length
$goOverPos = $sentenceStart + $goOverOffset

# Drop a zero-length bookmark exactly where "go over" will begin once
# "briefly " is removed. Doing this *before* the text edit forces the
# run to split there, so the later in-place delete only re-flows the
# left-hand half of the run (leaving the following " We assumed..."
# run untouched).
$bmRange = $d.Range($goOverPos, $goOverPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$briefPos = $sentenceStart + ("First, I want to ").Length
$rngBriefly = $d.Range($briefPos, $briefPos + ("briefly ").Length)
if ($rngBriefly.Text -ne "briefly ") {
    throw "unexpected text where 'briefly ' was expected: [$($rngBriefly.Text)]"
}
$rngBriefly.Text = ""

# ------------------------------------------------------------------
# 2) Append a new sentence to the final paragraph, as its own run.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$endRng = $lastPara.Range
$endRng.MoveEnd(1, -1) | Out-Null
$endRng.Collapse(0)

# Insert via a throw-away bookmark so the new sentence is materialized
# as its own run instead of being folded back into the preceding run
# (both happen to share the same, default run formatting).
$d.Bookmarks.Add("ZZZ_tmp_split", $endRng)
$insertRng = $d.Range($endRng.End, $endRng.End)
$newSentence = " Yi" + $apos + "s work has already yielded a sample TTL file, and it" + $apos + "s very representative of what the TTL files should look like."
$insertRng.InsertAfter($newSentence)
$d.Bookmarks("ZZZ_tmp_split").Delete()
